# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Both sheets list the same events; the row numbers differ by one offset
# starting at row 6 on "展览" (row 7 on "全部类型") because "全部类型"
# contains one extra row (row 5) that "展览" does not have.

$wb = $excel.ActiveWorkbook

# Row => New value, for worksheet "展览" (F2..F39)
$updatesSheet1 = @{
    2  = 7335
    6  = 207
    7  = 145
    10 = 69
    11 = 235
    12 = 22
    13 = 477
    14 = 33
    15 = 1885
    17 = 55
    18 = 3870
    19 = 31
    20 = 257
    21 = 87
    22 = 61
    24 = 43
    25 = 2569
    26 = 29
    27 = 342
    30 = 49
    31 = 11
    32 = 35
    34 = 3
    37 = 68
    38 = 1524
    39 = 186
}

# Row => New value, for worksheet "全部类型" (F2..F40)
$updatesSheet4 = @{
    2  = 7335
    5  = 9
    7  = 207
    8  = 145
    11 = 69
    12 = 235
    13 = 22
    14 = 477
    15 = 33
    16 = 1885
    18 = 55
    19 = 3870
    20 = 31
    21 = 257
    22 = 87
    23 = 61
    25 = 43
    26 = 2569
    27 = 29
    28 = 342
    31 = 49
    32 = 11
    33 = 35
    35 = 3
    38 = 68
    39 = 1524
    40 = 186
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updatesSheet1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updatesSheet1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesSheet4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updatesSheet4[$row]
}
